# New crime data collected
# This script updates the 104th Precinct CompStat weekly report:
#  - bumps the report Volume/Number and the "Week Covering" date range
#  - refreshes the weekly crime-complaint statistics table (rows 14-29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Volume/Number and reporting week dates ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Numeric value updates in the statistics table ---
$ws.Range("N14").Value = -89.473684210526
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 21
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 10.526315789473
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 213
$ws.Range("J16").Value = 135
$ws.Range("K16").Value = 57.777777777777
$ws.Range("L16").Value = 21.714285714285
$ws.Range("M16").Value = -17.441860465116
$ws.Range("N16").Value = -78.419452887538
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 295
$ws.Range("J17").Value = 248
$ws.Range("K17").Value = 18.951612903225
$ws.Range("L17").Value = 3.873239436619
$ws.Range("M17").Value = 36.574074074074
$ws.Range("N17").Value = -0.673400673400
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("I18").Value = 262
$ws.Range("J18").Value = 234
$ws.Range("K18").Value = 11.965811965812
$ws.Range("L18").Value = -15.210355987055
$ws.Range("M18").Value = -39.63133640553
$ws.Range("N18").Value = -86.361270171785
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -55
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -46.835443037974
$ws.Range("I19").Value = 619
$ws.Range("J19").Value = 558
$ws.Range("K19").Value = 10.931899641577
$ws.Range("L19").Value = 6.540447504302
$ws.Range("M19").Value = 47.380952380952
$ws.Range("N19").Value = -5.640243902439
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 27
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 306
$ws.Range("J20").Value = 207
$ws.Range("K20").Value = 47.826086956521
$ws.Range("L20").Value = 71.910112359550
$ws.Range("M20").Value = -13.314447592068
$ws.Range("N20").Value = -90.808050465605
$ws.Range("C21").Value = 41
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 5.128205128205
$ws.Range("F21").Value = 140
$ws.Range("G21").Value = 159
$ws.Range("H21").Value = -11.949685534591
$ws.Range("I21").Value = 1718
$ws.Range("J21").Value = 1399
$ws.Range("K21").Value = 22.802001429592
$ws.Range("L21").Value = 11.486048020765
$ws.Range("M21").Value = 0.821596244131
$ws.Range("N21").Value = -76.237897648686
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = 77.777777777777
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = 2.941176470588
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 147
$ws.Range("H24").Value = -8.843537414965
$ws.Range("I24").Value = 1447
$ws.Range("J24").Value = 1345
$ws.Range("K24").Value = 7.583643122676
$ws.Range("L24").Value = 1.188811188811
$ws.Range("M24").Value = 40.621963070942
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -5.882352941176
$ws.Range("I25").Value = 513
$ws.Range("J25").Value = 480
$ws.Range("K25").Value = 6.875
$ws.Range("L25").Value = 13.245033112582
$ws.Range("M25").Value = -23.774145616641
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 35
$ws.Range("K26").Value = 66.666666666666
$ws.Range("L26").Value = 29.629629629629
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 66
$ws.Range("K27").Value = 29.411764705882
$ws.Range("L27").Value = 43.478260869565
$ws.Range("N28").Value = -67.741935483871
$ws.Range("N29").Value = -67.857142857142

# --- Cells that become the literal text "0" (distinct from the numeric 0) ---
# NumberFormat is switched to Text ("@") first so the digit string is stored
# as text instead of being re-interpreted as a number.
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"

# --- Cells that become the "N/A" placeholder text ---
$ws.Range("E15").Value = "***.*"
$ws.Range("E26").Value = "***.*"
$ws.Range("E27").Value = "***.*"
